# Daily BP terminal gate pricing refresh:
# - Shift the "most recent" (45986) price rows down into the "previous"
#   (45983) slot, which they replace.
# - Populate the new top (45987) rows with the newly published prices.
# Only the data cells change; row/column formatting, styles, and the
# header/label rows are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 45987
$ws.Range("D8").Value = 171.96
$ws.Range("E8").Value = 160.74
$ws.Range("F8").Value = 170.74
$ws.Range("G8").Value = 160.91
$ws.Range("A9").Value = 45987
$ws.Range("D9").Value = 171.96
$ws.Range("E9").Value = 160.74
$ws.Range("F9").Value = 170.74
$ws.Range("G9").Value = 160.91
$ws.Range("A10").Value = 45987
$ws.Range("D10").Value = 173.71
$ws.Range("E10").Value = 163.03
$ws.Range("F10").Value = 173.03
$ws.Range("G10").Value = 163.52
$ws.Range("A11").Value = 45986
$ws.Range("D11").Value = 173.02
$ws.Range("E11").Value = 161.29
$ws.Range("F11").Value = 171.29
$ws.Range("G11").Value = 161.46
$ws.Range("A12").Value = 45986
$ws.Range("D12").Value = 173.02
$ws.Range("E12").Value = 161.29
$ws.Range("F12").Value = 171.29
$ws.Range("G12").Value = 161.46
$ws.Range("A13").Value = 45986
$ws.Range("D13").Value = 174.97
$ws.Range("E13").Value = 163.77
$ws.Range("F13").Value = 173.77
$ws.Range("G13").Value = 164.26
$ws.Range("A17").Value = 45987
$ws.Range("D17").Value = 176.98
$ws.Range("E17").Value = 165.5
$ws.Range("F17").Value = 175.5
$ws.Range("A18").Value = 45986
$ws.Range("D18").Value = 178.25
$ws.Range("E18").Value = 166.26
$ws.Range("F18").Value = 176.26
$ws.Range("A22").Value = 45987
$ws.Range("D22").Value = 173.09
$ws.Range("E22").Value = 162.19
$ws.Range("F22").Value = 171.79
$ws.Range("G22").Value = 163.48
$ws.Range("A23").Value = 45987
$ws.Range("D23").Value = 178.5
$ws.Range("E23").Value = 166.95
$ws.Range("F23").Value = 176.95
$ws.Range("A24").Value = 45987
$ws.Range("D24").Value = 178.28
$ws.Range("E24").Value = 167.23
$ws.Range("F24").Value = 177.23
$ws.Range("A25").Value = 45987
$ws.Range("D25").Value = 179.11
$ws.Range("E25").Value = 166.63
$ws.Range("F25").Value = 176.63
$ws.Range("G25").Value = 166.67
$ws.Range("A26").Value = 45987
$ws.Range("D26").Value = 177.8
$ws.Range("E26").Value = 168.08
$ws.Range("F26").Value = 178.08
$ws.Range("A27").Value = 45986
$ws.Range("D27").Value = 174.14
$ws.Range("E27").Value = 162.83
$ws.Range("F27").Value = 172.43
$ws.Range("G27").Value = 164.11
$ws.Range("A28").Value = 45986
$ws.Range("D28").Value = 179.76
$ws.Range("E28").Value = 167.59
$ws.Range("F28").Value = 177.59
$ws.Range("A29").Value = 45986
$ws.Range("D29").Value = 179.55
$ws.Range("E29").Value = 167.97
$ws.Range("F29").Value = 177.97
$ws.Range("A30").Value = 45986
$ws.Range("D30").Value = 180.37
$ws.Range("E30").Value = 167.38
$ws.Range("F30").Value = 177.38
$ws.Range("G30").Value = 167.42
$ws.Range("A31").Value = 45986
$ws.Range("D31").Value = 179.06
$ws.Range("E31").Value = 168.83
$ws.Range("F31").Value = 178.83
$ws.Range("A35").Value = 45987
$ws.Range("D35").Value = 172.06
$ws.Range("E35").Value = 160.01
$ws.Range("F35").Value = 169.02
$ws.Range("A36").Value = 45986
$ws.Range("D36").Value = 173.32
$ws.Range("E36").Value = 160.76
$ws.Range("F36").Value = 169.76
$ws.Range("A40").Value = 45987
$ws.Range("D40").Value = 177.74
$ws.Range("E40").Value = 165.51
$ws.Range("F40").Value = 175.51
$ws.Range("A41").Value = 45987
$ws.Range("D41").Value = 177.45
$ws.Range("E41").Value = 165.93
$ws.Range("F41").Value = 175.93
$ws.Range("A42").Value = 45986
$ws.Range("D42").Value = 179
$ws.Range("E42").Value = 166.26
$ws.Range("F42").Value = 176.26
$ws.Range("A43").Value = 45986
$ws.Range("D43").Value = 178.7
$ws.Range("E43").Value = 166.68
$ws.Range("F43").Value = 176.68
$ws.Range("A47").Value = 45987
$ws.Range("D47").Value = 173.43
$ws.Range("E47").Value = 162.31
$ws.Range("F47").Value = 172.31
$ws.Range("A48").Value = 45987
$ws.Range("D48").Value = 173.38
$ws.Range("E48").Value = 162.45
$ws.Range("F48").Value = 172.45
$ws.Range("A49").Value = 45986
$ws.Range("D49").Value = 173.02
$ws.Range("E49").Value = 162.59
$ws.Range("F49").Value = 172.59
$ws.Range("A50").Value = 45986
$ws.Range("D50").Value = 172.98
$ws.Range("E50").Value = 162.74
$ws.Range("F50").Value = 172.74
$ws.Range("A54").Value = 45987
$ws.Range("D54").Value = 187.79
$ws.Range("E54").Value = 176.42
$ws.Range("F54").Value = 186.42
$ws.Range("A55").Value = 45987
$ws.Range("D55").Value = 175.69
$ws.Range("E55").Value = 173.12
$ws.Range("F55").Value = 183.12
$ws.Range("A56").Value = 45987
$ws.Range("D56").Value = 177.98
$ws.Range("A57").Value = 45987
$ws.Range("D57").Value = 177.4
$ws.Range("E57").Value = 167.39
$ws.Range("A58").Value = 45987
$ws.Range("D58").Value = 173.31
$ws.Range("E58").Value = 163.44
$ws.Range("F58").Value = 173.44
$ws.Range("A59").Value = 45987
$ws.Range("D59").Value = 179.93
$ws.Range("E59").Value = 174.33
$ws.Range("A60").Value = 45986
$ws.Range("D60").Value = 189.05
$ws.Range("E60").Value = 177.14
$ws.Range("F60").Value = 187.14
$ws.Range("A61").Value = 45986
$ws.Range("D61").Value = 176.95
$ws.Range("E61").Value = 173.76
$ws.Range("F61").Value = 183.76
$ws.Range("A62").Value = 45986
$ws.Range("D62").Value = 179.13
$ws.Range("A63").Value = 45986
$ws.Range("D63").Value = 178.56
$ws.Range("E63").Value = 168.03
$ws.Range("A64").Value = 45986
$ws.Range("D64").Value = 174.47
$ws.Range("E64").Value = 164.08
$ws.Range("F64").Value = 174.08
$ws.Range("A65").Value = 45986
$ws.Range("D65").Value = 181.2
$ws.Range("E65").Value = 175.06
